# "update for 2023 classroom"
# - Remove student "Joohu" from the list (row 9), shifting the names
#   below it up by one row.
# - Simplify the last student's display name from "Yonghan (Addie)" to
#   "Addie" (now landing on row 14 after the shift).
# - Clear the now-unused last row (15) in column A.
# - Update the sheet selection to the block of names that moved (A9:A14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value  = "Cameron"
$ws.Range("A10").Value = "Ian"
$ws.Range("A11").Value = "John"
$ws.Range("A12").Value = "Nima"
$ws.Range("A13").Value = "Sai"
$ws.Range("A14").Value = "Addie"
$ws.Range("A15").ClearContents()

$ws.Range("A9:A14").Select() | Out-Null
